# "Generate Report for Handoff"
#
# The aa0ab8e1-ea6a-411d-9970-d3af90ea867b.* row is removed from every
# sheet (it finished handback, so it's gone from this handoff report),
# the 9675f9d4-...-e4d43d867617 row's status flips from the old
# "Handed back: in sync with en-US" text to "Ready for handoff", and its
# handoff datetimes are refreshed.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet "Overview": A1:C4 -> A1:C3
# ---------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Overview")

$ws1.Hyperlinks.Delete()
$ws1.Rows.Item(3).Delete()

$ws1.Range("B2").Value = "Ready for handoff"
$ws1.Range("C2").Value = "Ready for handoff"

$ws1.Hyperlinks.Add($ws1.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/004f581602244b6fda38598a580fdb4ac2b7e7e0/e2e/9675f9d4-d77a-4429-af2a-e4d43d867617.md", "", "", "9675f9d4-d77a-4429-af2a-e4d43d867617.md")
$ws1.Hyperlinks.Add($ws1.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/004f581602244b6fda38598a580fdb4ac2b7e7e0/.localization-config", "", "", ".localization-config")

# ---------------------------------------------------------------
# Sheet "zh-cn": A1:I4 -> A1:I3
# ---------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("zh-cn")

$ws2.Hyperlinks.Delete()
$ws2.Rows.Item(3).Delete()

$ws2.Range("B2").Value = "Ready for handoff"
$ws2.Range("D2").Value = "2016-01-25 07:56:38"
$ws2.Range("D2").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws2.Range("D3").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$ws2.Hyperlinks.Add($ws2.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/004f581602244b6fda38598a580fdb4ac2b7e7e0/e2e/9675f9d4-d77a-4429-af2a-e4d43d867617.md", "", "", "9675f9d4-d77a-4429-af2a-e4d43d867617.md")
$ws2.Hyperlinks.Add($ws2.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/5b74838cb6509a61d4af5303ef4c1e6ee72495d9/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/yuwzho/9675f9d4-d77a-4429-af2a-e4d43d867617.d9be021a9c0fdb26d74d342892fee51c1d062707.zh-cn.xlf", "", "", "9675f9d4-d77a-4429-af2a-e4d43d867617.d9be021a9c0fdb26d74d342892fee51c1d062707.zh-cn.xlf")
$ws2.Hyperlinks.Add($ws2.Range("E2"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/ae09dd006d92233dd5c7be7751d6e4c99c05a945/e2e/9675f9d4-d77a-4429-af2a-e4d43d867617.md", "", "", "9675f9d4-d77a-4429-af2a-e4d43d867617.md")
$ws2.Hyperlinks.Add($ws2.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/e453890af13f5f412d696bb1278fae5d7831d844/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/yuwzho/9675f9d4-d77a-4429-af2a-e4d43d867617.d9be021a9c0fdb26d74d342892fee51c1d062707.zh-cn.xlf", "", "", "9675f9d4-d77a-4429-af2a-e4d43d867617.d9be021a9c0fdb26d74d342892fee51c1d062707.zh-cn.xlf")
$ws2.Hyperlinks.Add($ws2.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/004f581602244b6fda38598a580fdb4ac2b7e7e0/.localization-config", "", "", ".localization-config")

# ---------------------------------------------------------------
# Sheet "de-de": A1:I4 -> A1:I3
# ---------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("de-de")

$ws3.Hyperlinks.Delete()
$ws3.Rows.Item(3).Delete()

$ws3.Range("B2").Value = "Ready for handoff"
$ws3.Range("D2").Value = "2016-01-25 07:56:50"
$ws3.Range("D2").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws3.Range("D3").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$ws3.Hyperlinks.Add($ws3.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/004f581602244b6fda38598a580fdb4ac2b7e7e0/e2e/9675f9d4-d77a-4429-af2a-e4d43d867617.md", "", "", "9675f9d4-d77a-4429-af2a-e4d43d867617.md")
$ws3.Hyperlinks.Add($ws3.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/209fa3e3c8d2c5c035715b037c878e0300bc8789/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/yuwzho/9675f9d4-d77a-4429-af2a-e4d43d867617.d9be021a9c0fdb26d74d342892fee51c1d062707.de-de.xlf", "", "", "9675f9d4-d77a-4429-af2a-e4d43d867617.d9be021a9c0fdb26d74d342892fee51c1d062707.de-de.xlf")
$ws3.Hyperlinks.Add($ws3.Range("E2"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/c9c2b3553882f449792ec8fb9d49087ec58a4787/e2e/9675f9d4-d77a-4429-af2a-e4d43d867617.md", "", "", "9675f9d4-d77a-4429-af2a-e4d43d867617.md")
$ws3.Hyperlinks.Add($ws3.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/537f36c7606cfa75b0d33035dee9204520e0aa26/ol-handback/OpenLocalizationTestOrg/oltest.de-de/yuwzho/9675f9d4-d77a-4429-af2a-e4d43d867617.d9be021a9c0fdb26d74d342892fee51c1d062707.de-de.xlf", "", "", "9675f9d4-d77a-4429-af2a-e4d43d867617.d9be021a9c0fdb26d74d342892fee51c1d062707.de-de.xlf")
$ws3.Hyperlinks.Add($ws3.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/004f581602244b6fda38598a580fdb4ac2b7e7e0/.localization-config", "", "", ".localization-config")
